# Base clientes.xlsx - agrega 4 nuevos clientes a la tabla "Tabla1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Agregar los nuevos registros de clientes ---
$ws.Range("A16").Value = 10240057
$ws.Range("B16").Value = "MARIA GUADALUPE RUIZ ESTRADA"

$ws.Range("A17").Value = 10174274
$ws.Range("B17").Value = "AGRICOLA EL MORON SA DE CV"

$ws.Range("A18").Value = 500231
$ws.Range("B18").Value = "INSECTICIDAS HERBICIDAS Y SEMILLAS RIDA"

$ws.Range("A19").Value = 10181721
$ws.Range("B19").Value = "DUNE COMPANY MEXICALI"

# --- Expandir la tabla Tabla1 para incluir las nuevas filas ---
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("A1:B19"))

# --- Ajustar el ancho de la columna B para mostrar el texto mas largo ---
$ws.Columns("B:B").ColumnWidth = 38.21875

# --- Actualizar la celda seleccionada (ultima posicion usada al editar) ---
$ws.Range("B14").Select() | Out-Null
